$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31
$ws.Range("G31").Value = 2.63
$ws.Range("H31").Value = 2.9
$ws.Range("I31").Value = 3
$ws.Range("J31").Value = 3.4
$ws.Range("K31").Value = 1.91
$ws.Range("O31").Value = 1.5
$ws.Range("P31").Value = 2.5
$ws.Range("Q31").Value = 2.6
$ws.Range("R31").Value = 1.48
$ws.Range("S31").Value = 4.2
$ws.Range("T31").Value = 1.23
$ws.Range("AC31").Value = 11
$ws.Range("AD31").Value = 26
$ws.Range("AE31").Value = 26
$ws.Range("AL31").Value = 7
$ws.Range("AO31").Value = 29
$ws.Range("AR31").Value = 1.98
$ws.Range("AS31").Value = 1.88

# Row 36
$ws.Range("G36").Value = 2.05
$ws.Range("I36").Value = 4.33
$ws.Range("J36").Value = 3
$ws.Range("O36").Value = 1.83
$ws.Range("P36").Value = 1.83
$ws.Range("AA36").Value = 4.5
$ws.Range("AD36").Value = 19
$ws.Range("AL36").Value = 7
$ws.Range("AM36").Value = 19
$ws.Range("AQ36").Value = 67

# Row 55
$ws.Range("G55").Value = 1.53
$ws.Range("H55").Value = 4.75
$ws.Range("I55").Value = 5
$ws.Range("J55").Value = 2.1
$ws.Range("K55").Value = 2.38
$ws.Range("O55").Value = 1.25
$ws.Range("P55").Value = 3.75
$ws.Range("Q55").Value = 1.75
$ws.Range("R55").Value = 2.05
$ws.Range("S55").Value = 2.23
$ws.Range("U55").Value = 3
$ws.Range("V55").Value = 1.36
$ws.Range("AD55").Value = 11
$ws.Range("AG55").Value = 13
$ws.Range("AH55").Value = 9
$ws.Range("AI55").Value = 19
$ws.Range("AL55").Value = 13
$ws.Range("AM55").Value = 26
$ws.Range("AN55").Value = 15
$ws.Range("AO55").Value = 51

# Row 57
$ws.Range("H57").Value = 2.88
$ws.Range("I57").Value = 3.6
$ws.Range("M57").Value = 1.13
$ws.Range("N57").Value = 6
$ws.Range("O57").Value = 1.53
$ws.Range("P57").Value = 2.38
$ws.Range("Q57").Value = 2.7
$ws.Range("R57").Value = 1.44
$ws.Range("W57").Value = 1.62
$ws.Range("X57").Value = 2.2
$ws.Range("Y57").Value = 2.2
$ws.Range("Z57").Value = 1.62
$ws.Range("AH57").Value = 6
$ws.Range("AI57").Value = 19
$ws.Range("AM57").Value = 15
$ws.Range("AQ57").Value = 51
$ws.Range("AR57").Value = 2.03
$ws.Range("AS57").Value = 1.83

# Row 68
$ws.Range("G68").Value = 2.1
$ws.Range("H68").Value = 3.1
$ws.Range("I68").Value = 3.7
$ws.Range("J68").Value = 2.88
$ws.Range("K68").Value = 2
$ws.Range("O68").Value = 1.44
$ws.Range("P68").Value = 2.63
$ws.Range("Q68").Value = 2.35
$ws.Range("R68").Value = 1.57
$ws.Range("U68").Value = 4.5
$ws.Range("V68").Value = 1.18
$ws.Range("W68").Value = 1.53
$ws.Range("X68").Value = 2.38
$ws.Range("Y68").Value = 2
$ws.Range("Z68").Value = 1.73
$ws.Range("AA68").Value = 6
$ws.Range("AC68").Value = 9.5
$ws.Range("AD68").Value = 19
$ws.Range("AG68").Value = 7
$ws.Range("AJ68").Value = 67
$ws.Range("AL68").Value = 9
$ws.Range("AM68").Value = 17
$ws.Range("AR68").Value = 1.78
$ws.Range("AS68").Value = 2.03

# Row 140
$ws.Range("I140").Value = 3.9
$ws.Range("M140").Value = 1.07
$ws.Range("N140").Value = 9
$ws.Range("O140").Value = 1.4
$ws.Range("P140").Value = 2.75
$ws.Range("Q140").Value = 2.25
$ws.Range("R140").Value = 1.62
$ws.Range("U140").Value = 4.33
$ws.Range("V140").Value = 1.2
$ws.Range("Y140").Value = 2
$ws.Range("Z140").Value = 1.75
$ws.Range("AB140").Value = 8.5
$ws.Range("AE140").Value = 19

# Row 141
$ws.Range("G141").Value = 2.4
$ws.Range("I141").Value = 3.1
$ws.Range("J141").Value = 3.1
$ws.Range("K141").Value = 2.1
$ws.Range("U141").Value = 3.5
$ws.Range("V141").Value = 1.29
$ws.Range("AC141").Value = 9.5
$ws.Range("AL141").Value = 9.5
$ws.Range("AN141").Value = 12
$ws.Range("AO141").Value = 34
$ws.Range("AP141").Value = 26

# Row 142
$ws.Range("G142").Value = 2.15
$ws.Range("I142").Value = 3.1
$ws.Range("J142").Value = 2.77
$ws.Range("L142").Value = 3.65
$ws.Range("Q142").Value = 1.85
$ws.Range("U142").Value = 2.95
$ws.Range("X142").Value = 2.52
$ws.Range("Z142").Value = 1.93
$ws.Range("AA142").Value = 7.9
$ws.Range("AB142").Value = 10.75
$ws.Range("AC142").Value = 8.75
$ws.Range("AD142").Value = 21
$ws.Range("AE142").Value = 17
$ws.Range("AF142").Value = 27
$ws.Range("AL142").Value = 9.5
$ws.Range("AM142").Value = 16
$ws.Range("AN142").Value = 11
$ws.Range("AO142").Value = 40
$ws.Range("AP142").Value = 27
$ws.Range("AQ142").Value = 35

# Row 148
$ws.Range("H148").Value = 4.33
$ws.Range("I148").Value = 9
$ws.Range("Q148").Value = 1.95
$ws.Range("R148").Value = 1.85
$ws.Range("U148").Value = 3.4
$ws.Range("V148").Value = 1.3
$ws.Range("AL148").Value = 19

# Row 152
$ws.Range("J152").Value = 1.67
$ws.Range("L152").Value = 11
$ws.Range("Q152").Value = 1.83
$ws.Range("R152").Value = 2.03
$ws.Range("U152").Value = 3
$ws.Range("V152").Value = 1.36
$ws.Range("AB152").Value = 5.5
$ws.Range("AD152").Value = 7
$ws.Range("AJ152").Value = 101
$ws.Range("AL152").Value = 26
$ws.Range("AP152").Value = 101
$ws.Range("AR152").Value = 1.41
$ws.Range("AS152").Value = 2.95

# Row 207
$ws.Range("K207").Value = 1.95
$ws.Range("O207").Value = 1.44
$ws.Range("P207").Value = 2.63
$ws.Range("Q207").Value = 2.4
$ws.Range("R207").Value = 1.53
$ws.Range("AA207").Value = 6.5
$ws.Range("AH207").Value = 6
$ws.Range("AR207").Value = 1.88
$ws.Range("AS207").Value = 1.98

# Row 210
$ws.Range("G210").Value = 2.47
$ws.Range("H210").Value = 3.1
$ws.Range("I210").Value = 2.75
$ws.Range("J210").Value = 3.2
$ws.Range("K210").Value = 1.98
$ws.Range("L210").Value = 3.35
$ws.Range("U210").Value = 3.3
$ws.Range("V210").Value = 1.24
$ws.Range("X210").Value = 2.37
$ws.Range("Y210").Value = 1.78
$ws.Range("Z210").Value = 1.83
$ws.Range("AA210").Value = 7.3
$ws.Range("AB210").Value = 11.5
$ws.Range("AC210").Value = 9.75
$ws.Range("AD210").Value = 27
$ws.Range("AE210").Value = 22
$ws.Range("AG210").Value = 8.25
$ws.Range("AH210").Value = 6
$ws.Range("AI210").Value = 14.5
$ws.Range("AJ210").Value = 75
$ws.Range("AK210").Value = 600
$ws.Range("AL210").Value = 8.25
$ws.Range("AN210").Value = 10
$ws.Range("AP210").Value = 24
$ws.Range("AQ210").Value = 35
